$d = $word.ActiveDocument

$replacements = @(
    @("527×6=", "949×6="),
    @("485×9=", "674×6="),
    @("658×9=", "126×3="),
    @("625×9=", "884×6="),
    @("586×4=", "842×4="),
    @("380×5=", "277×5="),
    @("509×7=", "962×9="),
    @("614×6=", "137×7="),
    @("143×5=", "622×4="),
    @("970×3=", "814×7="),
    @("500×7=", "160×4="),
    @("264×6=", "751×5="),
    @("449×5=", "197×3="),
    @("921×3=", "551×7="),
    @("118×7=", "236×4="),
    @("295×2=", "750×3="),
    @("855×6=", "486×7="),
    @("948×6=", "979×2="),
    @("456×3=", "984×8="),
    @("929×7=", "534×8="),
    @("584×8=", "811×7="),
    @("164×2=", "830×6="),
    @("393×8=", "835×9="),
    @("258×8=", "261×4="),
    @("676×6=", "396×5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
